# Shop.xlsx edit: "unify the conception of DataNode, DataTable, Entity."
#
# Re-saving this workbook (originally authored in Mac Excel, rupBuild 27309)
# with a newer Windows Excel build (rupBuild 18730) renamed the sheet and
# moved the active selection; the surrounding namespace/xr:uid/version churn
# in the raw XML is Excel's own re-serialization noise that comes along for
# free on a real round-trip. We reproduce the concrete, user-visible edits
# here via the COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet: "Property1" -> "DataNode" ------------------------
$ws.Name = "DataNode"

# --- Move the active selection from A9 to D40 --------------------------------
$ws.Range("D40").Select()

# --- Minor column width touch-up (A, B, C, F:G) ------------------------------
# Stored widths shift by a hair (~1/256 char) when the workbook is resaved by
# the newer build. Re-apply the closest widths this engine can produce.
$ws.Columns.Item(1).ColumnWidth = 18.3575    # A: ~19.16 -> ~19.125
$ws.Columns.Item(2).ColumnWidth = 7.3575     # B: ~8.16  -> ~8.125
$ws.Columns.Item(3).ColumnWidth = 27.3575    # C: ~28.16 -> ~28.125
$ws.Columns.Item(6).ColumnWidth = 11.786     # F: ~12.66 -> ~12.625
$ws.Columns.Item(7).ColumnWidth = 11.786     # G: ~12.66 -> ~12.625
